# Fruta / hortaliza, semanal
# Insert a new weekly record at row 184 (pushing the existing rows 184-203
# down to 185-204) in the "Femacal de La Calera - Poroto verde" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 184, shifting rows 184:203 down to 185:204.
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new observation.
$ws.Range("A184").Value = 3
$ws.Range("B184").Value = "Femacal de La Calera"
$ws.Range("C184").Value = "Coquimbo"
$ws.Range("D184").Value = 44449
$ws.Range("E184").Value = 5
$ws.Range("F184").Value = 100112031
$ws.Range("G184").Value = "Poroto verde"
$ws.Range("H184").Value = "Magnum"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 73
$ws.Range("K184").Value = 32000
$ws.Range("L184").Value = 33000
$ws.Range("M184").Value = 32479
$ws.Range("N184").Value = "`$/malla 25 kilos"
$ws.Range("O184").Value = "Región de Arica y Parinacota"
$ws.Range("P184").Value = 1299
$ws.Range("Q184").Value = 25
$ws.Range("R184").Value = "Hortaliza"
